# ToDoItems.xlsx fixture update:
# The todo-item rows 9-13 (penpal/blog/brown-bag/conference/Isis-release
# "professional" sample rows) are no longer part of the fixture data, so
# their category/subcategory/description/cost contents are cleared out,
# leaving just the pre-existing row style on column A. The remaining rows
# (and their shared-string references) stay as-is; clearing the now-unused
# strings causes the surviving shared strings to be renumbered.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9/10 column A previously held a (redundant) category value too -
# clear it along with the rest of the now-removed sample rows.
$ws.Range("A9:A10").ClearContents()
$ws.Range("B9:D13").ClearContents()

# Move the active selection to A8 (previously D2:D4).
[void]$ws.Activate()
$ws.Range("A8").Select() | Out-Null
